$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 101105.3
$ws.Range("I41").Value = 75.40000000000001
$ws.Range("J41").Value = 202135.2
$ws.Range("K41").Value = 75.40000000000001
$ws.Range("L41").Value = 202135.2
$ws.Range("M41").Value = 364.6
$ws.Range("N41").Value = -203015.2
$ws.Range("H106").Value = 3171.2856
$ws.Range("J106").Value = 4933
$ws.Range("L106").Value = 4933
$ws.Range("N106").Value = -6195
$ws.Range("H132").Value = 55933.934
$ws.Range("I132").Value = 69032.88
$ws.Range("K132").Value = 207098.64
$ws.Range("M132").Value = -204568.64
$ws.Range("H137").Value = 2026
$ws.Range("I137").Value = 1710.6857
$ws.Range("K137").Value = 5132.0571
$ws.Range("M137").Value = -2582.0571
$ws.Range("H138").Value = 1703.7609
$ws.Range("I138").Value = 1100.909
$ws.Range("J138").Value = 3234.077
$ws.Range("K138").Value = 3302.727
$ws.Range("L138").Value = 9702.231
$ws.Range("M138").Value = 1837.273
$ws.Range("N138").Value = -19982.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10871536
$ws.Range("I32").Value = 11629520
$ws.Range("K32").Value = 11629520
$ws.Range("M32").Value = -11629233
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H61").Value = 8763.25
$ws.Range("I61").Value = 13052.667
$ws.Range("J61").Value = 5253.727
$ws.Range("K61").Value = 13052.667
$ws.Range("L61").Value = 5253.727
$ws.Range("M61").Value = -12840.667
$ws.Range("N61").Value = -5677.727
$ws.Range("H92").Value = 87033.336
$ws.Range("J92").Value = 87033.336
$ws.Range("L92").Value = 87033.336
$ws.Range("N92").Value = -92025.336
$ws.Range("H102").Value = 14253.156
$ws.Range("I102").Value = 15396.621
$ws.Range("J102").Value = 3199.6667
$ws.Range("K102").Value = 15396.621
$ws.Range("L102").Value = 3199.6667
$ws.Range("M102").Value = -13774.621
$ws.Range("N102").Value = -6443.6667
$ws.Range("H136").Value = 8763.25
$ws.Range("I136").Value = 13052.667
$ws.Range("J136").Value = 5253.727
$ws.Range("K136").Value = 39158.001
$ws.Range("L136").Value = 15761.181
$ws.Range("M136").Value = -36608.001
$ws.Range("N136").Value = -20861.181
$ws.Range("H141").Value = 82666
$ws.Range("J141").Value = 82666
$ws.Range("L141").Value = 82666
$ws.Range("N141").Value = -93026

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 22225488
$ws.Range("I107").Value = 3635.875
$ws.Range("K107").Value = 3635.875
$ws.Range("M107").Value = -1715.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 58022.89
$ws.Range("I31").Value = 76341.39999999999
$ws.Range("K31").Value = 76341.39999999999
$ws.Range("M31").Value = -76046.39999999999
$ws.Range("H34").Value = 58022.89
$ws.Range("I34").Value = 76341.39999999999
$ws.Range("K34").Value = 76341.39999999999
$ws.Range("M34").Value = -76139.39999999999
$ws.Range("H107").Value = 646.7857
$ws.Range("I107").Value = 355.1111
$ws.Range("K107").Value = 355.1111
$ws.Range("M107").Value = 1564.8889
$ws.Range("H122").Value = 2151
$ws.Range("I122").Value = 1454.8182
$ws.Range("J122").Value = 3682.6
$ws.Range("K122").Value = 4364.4546
$ws.Range("L122").Value = 11047.8
$ws.Range("M122").Value = -1914.4546
$ws.Range("N122").Value = -15947.8
$ws.Range("H129").Value = 47482.5
$ws.Range("I129").Value = 44980
$ws.Range("K129").Value = 44980
$ws.Range("M129").Value = -39980
$ws.Range("H132").Value = 6414606.5
$ws.Range("I132").Value = 4967.968
$ws.Range("J132").Value = 31251954
$ws.Range("K132").Value = 14903.904
$ws.Range("L132").Value = 93755862
$ws.Range("M132").Value = -12373.904
$ws.Range("N132").Value = -93760922
$ws.Range("H134").Value = 14609.923
$ws.Range("I134").Value = 22970
$ws.Range("J134").Value = 3209.818
$ws.Range("K134").Value = 68910
$ws.Range("L134").Value = 9629.454000000002
$ws.Range("M134").Value = -66375
$ws.Range("N134").Value = -14699.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 381.2
$ws.Range("J23").Value = 366.7143
$ws.Range("L23").Value = 1100.1429
$ws.Range("N23").Value = -1570.1429
$ws.Range("H34").Value = 920.8889
$ws.Range("I34").Value = 912.5714
$ws.Range("J34").Value = 950
$ws.Range("K34").Value = 2737.7142
$ws.Range("L34").Value = 2850
$ws.Range("M34").Value = -2653.7142
$ws.Range("N34").Value = -3018
$ws.Range("H39").Value = 4580.2
$ws.Range("J39").Value = 4975.25
$ws.Range("L39").Value = 14925.75
$ws.Range("N39").Value = -15513.75
$ws.Range("H107").Value = 492.6154
$ws.Range("I107").Value = 373.33334
$ws.Range("K107").Value = 1120.00002
$ws.Range("M107").Value = 799.9999800000001
$ws.Range("H112").Value = 3746
$ws.Range("I112").Value = 1277
$ws.Range("K112").Value = 3831
$ws.Range("M112").Value = -2723
$ws.Range("H117").Value = 4600.55
$ws.Range("I117").Value = 1690.75
$ws.Range("J117").Value = 6540.4165
$ws.Range("K117").Value = 5072.25
$ws.Range("L117").Value = 19621.2495
$ws.Range("M117").Value = -1630.25
$ws.Range("N117").Value = -26505.2495
$ws.Range("H118").Value = 3098.5
$ws.Range("I118").Value = 2099.5
$ws.Range("J118").Value = 4097.5
$ws.Range("K118").Value = 6298.5
$ws.Range("L118").Value = 12292.5
$ws.Range("M118").Value = -5055.5
$ws.Range("N118").Value = -14778.5
$ws.Range("H137").Value = 6511.067
$ws.Range("I137").Value = 1947
$ws.Range("J137").Value = 8170.727
$ws.Range("K137").Value = 5841
$ws.Range("L137").Value = 24512.181
$ws.Range("M137").Value = -741
$ws.Range("N137").Value = -34712.181
$ws.Range("H139").Value = 2694.4348
$ws.Range("I139").Value = 726.2
$ws.Range("K139").Value = 2178.6
$ws.Range("M139").Value = 2961.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 11773.818
$ws.Range("I102").Value = 14438.125
$ws.Range("K102").Value = 14438.125
$ws.Range("M102").Value = -12816.125
$ws.Range("H122").Value = 4390.96
$ws.Range("I122").Value = 2037
$ws.Range("J122").Value = 16749.25
$ws.Range("K122").Value = 6111
$ws.Range("L122").Value = 50247.75
$ws.Range("M122").Value = -3661
$ws.Range("N122").Value = -55147.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1451.1428
$ws.Range("I22").Value = 1341.5
$ws.Range("J22").Value = 1636.6923
$ws.Range("K22").Value = 1341.5
$ws.Range("L22").Value = 1636.6923
$ws.Range("M22").Value = -1046.5
$ws.Range("N22").Value = -2226.6923
$ws.Range("H27").Value = 1451.1428
$ws.Range("I27").Value = 1341.5
$ws.Range("J27").Value = 1636.6923
$ws.Range("K27").Value = 1341.5
$ws.Range("L27").Value = 1636.6923
$ws.Range("M27").Value = -1234.5
$ws.Range("N27").Value = -1850.6923
$ws.Range("H46").Value = 1001
$ws.Range("I46").Value = 1001
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1001
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -813
$ws.Range("N46").ClearContents()
$ws.Range("H132").Value = 788163.9
$ws.Range("I132").Value = 1117650.6
$ws.Range("K132").Value = 3352951.8
$ws.Range("M132").Value = -3350421.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2358.9375
$ws.Range("I122").Value = 2124.5715
$ws.Range("K122").Value = 6373.7145
$ws.Range("M122").Value = -3923.7145
$ws.Range("H132").Value = 463833.2
$ws.Range("I132").Value = 630795.3
$ws.Range("K132").Value = 1892385.9
$ws.Range("M132").Value = -1889855.9
